# work on data processor
# Swap the contents of rows 2 and 3 (both the numeric "MatchValue" column A
# and the "ExampleData" label column B), and move the active selection from
# A3 to C4, matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 3
$ws.Range("A2").Value = 14.392
$ws.Range("B2").Value = "PID1_B"
$ws.Range("A3").Value = 12.321
$ws.Range("B3").Value = "PID1_A"

# Move the selection to C4
[void]$ws.Range("C4").Select()
